# Parametrized test: drop the free-text "test description" / duplicated
# "excel1@gmx" / "passwordexcel" helper cells from the manual-test sheet,
# keeping just the E-Mail/Passwort header row and the (now re-used) login
# value in A2.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# C2 held "invalid email FORMAT in the manual test" - no longer needed.
$ws.Range("C2").ClearContents()

# B3 ("passwordexcel") and C3 ("invalid email FORMAT in the manual test")
# are no longer needed either; A3 stays as-is (blank, styled).
$ws.Range("B3:C3").ClearContents()

# Move the saved selection like the authored workbook.
$ws.Range("G19").Select()
